$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell $ws "D2" "63.604.80"
Set-TextCell $ws "E2" "  -3.48%  "
Set-TextCell $ws "D3" "2.603.40"
Set-TextCell $ws "E3" "  -2.41%  "
Set-TextCell $ws "D4" "0.999"
Set-TextCell $ws "E4" "  -0.12%  "
Set-TextCell $ws "D5" "572.40"
Set-TextCell $ws "E5" "  -4.40%  "
Set-TextCell $ws "D6" "155.40"
Set-TextCell $ws "E6" "  -1.55%  "
Set-TextCell $ws "D7" "0.999"
Set-TextCell $ws "E7" "  -0.05%  "
Set-TextCell $ws "D8" "0.620"
Set-TextCell $ws "E8" "  -4.96%  "
Set-TextCell $ws "E9" "  -7.10%  "
Set-TextCell $ws "D10" "5.83"
Set-TextCell $ws "E10" "  -0.59%  "
Set-TextCell $ws "E11" "  -5.19%  "
Set-TextCell $ws "E12" "  -0.75%  "
Set-TextCell $ws "D13" "28.17"
Set-TextCell $ws "E13" "  -3.28%  "
Set-TextCell $ws "D14" "3.078.76"
Set-TextCell $ws "E14" "  -2.15%  "
Set-TextCell $ws "E15" "  -8.53%  "
Set-TextCell $ws "D16" "63.372.59"
Set-TextCell $ws "E16" "  -3.68%  "
Set-TextCell $ws "D17" "2.601.83"
Set-TextCell $ws "E17" "  -3.14%  "
Set-TextCell $ws "D18" "11.98"
Set-TextCell $ws "E18" "  -5.44%  "
Set-TextCell $ws "D19" "7.55"
Set-TextCell $ws "E19" "  +0.70%  "
Set-TextCell $ws "D20" "4.55"
Set-TextCell $ws "E20" "  -5.34%  "
Set-TextCell $ws "D21" "342.42"
Set-TextCell $ws "E21" "  -2.75%  "
Set-TextCell $ws "E22" "  +0.15%  "
Set-TextCell $ws "D23" "67.17"
Set-TextCell $ws "E23" "  -3.90%  "
Set-TextCell $ws "D24" "1.79"
Set-TextCell $ws "E24" "  -2.48%  "
Set-TextCell $ws "E25" "  -3.66%  "
Set-TextCell $ws "D26" "588.04"
Set-TextCell $ws "E26" "  +2.05%  "
Set-TextCell $ws "D27" "9.13"
Set-TextCell $ws "E27" "  -5.71%  "
Set-TextCell $ws "E28" "  -4.42%  "
Set-TextCell $ws "E29" "  +0.13%  "
Set-TextCell $ws "E30" "  -2.16%  "
Set-TextCell $ws "D31" "7.89"
Set-TextCell $ws "E31" "  -4.31%  "
Set-TextCell $ws "E32" "  -4.44%  "
Set-TextCell $ws "D33" "1.73"
Set-TextCell $ws "E33" "  -5.46%  "
Set-TextCell $ws "D34" "6.51"
Set-TextCell $ws "E34" "  -2.76%  "
Set-TextCell $ws "D35" "5.43"
Set-TextCell $ws "E35" "  -2.78%  "
Set-TextCell $ws "E36" "  -4.49%  "
Set-TextCell $ws "E37" "  -0.05%  "
Set-TextCell $ws "D38" "19.67"
Set-TextCell $ws "E38" "  -4.46%  "
Set-TextCell $ws "D39" "155.45"
Set-TextCell $ws "E39" "  +0.87%  "
Set-TextCell $ws "E40" "  -4.82%  "
Set-TextCell $ws "E41" "  +0.01%  "
Set-TextCell $ws "E42" "  +6.91%  "
Set-TextCell $ws "D43" "41.28"
Set-TextCell $ws "E43" "  -3.47%  "
Set-TextCell $ws "D44" "156.11"
Set-TextCell $ws "E44" "  -3.31%  "
Set-TextCell $ws "D45" "3.91"
Set-TextCell $ws "E45" "  -4.74%  "
Set-TextCell $ws "D46" "23.10"
Set-TextCell $ws "E46" "  +0.03%  "
Set-TextCell $ws "D47" "0.0587"
Set-TextCell $ws "E47" "  -5.17%  "
Set-TextCell $ws "D48" "0.628"
Set-TextCell $ws "E48" "  -2.55%  "
Set-TextCell $ws "E49" "  -1.77%  "
Set-TextCell $ws "D50" "0.0247"
Set-TextCell $ws "E50" "  -3.81%  "
Set-TextCell $ws "D51" "18.77"
Set-TextCell $ws "E51" "  -5.47%  "
